$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final target data (rows 2-10), reflecting removal of the original
# "peiper born 30 January 1915" row, an upward shift of the remaining
# rows, and two new rows appended for himler:Person (born / die).
$data = @(
    @("peiper:Person", "is_a",     "soldier"),
    @("peiper:Person", "die",      "14 July 1976"),
    @("peiper:Person", "is_a",     "nazi member"),
    @("peiper:Person", "is_a",     "German"),
    @("peiper:Person", "is_a",     "officer"),
    @("peiper:Person", "born_in",  "Wilmersdorf, Berlin, Germany"),
    @("himler:Person",  "know",    "peiper"),
    @("himler:Person",  "born",    "7 October 1900"),
    @("himler:Person",  "die",     "23 May 1945")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
